$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Phone / Price values for rows 6-24 (rows 1-5 and 25 are unchanged)
$data = @(
    @(6,  "Apple iPhone 6s Plus (Space Grey, 32 GB)", "₹34,900"),
    @(7,  "Apple iPhone 6s Plus (Silver, 32 GB)",      "₹34,900"),
    @(8,  "Apple iPhone XR (Black, 64 GB)",             "₹59,900"),
    @(9,  "Apple iPhone 6s (Gold, 32 GB)",               "₹29,799"),
    @(10, "Apple iPhone XR (Black, 128 GB)",             "₹64,900"),
    @(11, "Apple iPhone XR ((PRODUCT)RED, 64 GB)",       "₹59,490"),
    @(12, "Apple iPhone 6s (Space Grey, 32 GB)",         "₹29,590"),
    @(13, "Apple iPhone XR (White, 64 GB)",              "₹59,900"),
    @(14, "Apple iPhone XR (Blue, 64 GB)",                "₹59,900"),
    @(15, "Apple iPhone XR (White, 128 GB)",             "₹64,900"),
    @(16, "Apple iPhone XR (Coral, 128 GB)",             "₹64,900"),
    @(17, "Apple iPhone 7 (Silver, 128 GB)",              "₹48,999"),
    @(18, "Apple iPhone 6s Plus (Gold, 32 GB)",           "₹39,990"),
    @(19, "Apple iPhone XR (Blue, 128 GB)",               "₹64,900"),
    @(20, "Apple iPhone 7 (Rose Gold, 128 GB)",           "₹52,990"),
    @(21, "Apple iPhone 6s Plus (Silver, 16 GB)",         "₹39,990"),
    @(22, "Apple iPhone 7 (Silver, 32 GB)",               "₹39,900"),
    @(23, "Apple iPhone XR (Yellow, 256 GB)",             "₹74,900"),
    @(24, "Apple iPhone XR (Coral, 64 GB)",               "₹59,900")
)

foreach ($entry in $data) {
    $row = $entry[0]
    $phone = $entry[1]
    $price = $entry[2]
    $ws.Cells.Item($row, 1).Value = $phone
    $ws.Cells.Item($row, 2).Value = $price
}
